$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $excel.Windows.Count
$win = $excel.Windows.Item(1)
$win.ScrollRow = 42
Write-Host "done"
